$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The edited text lives in the speaker notes of slide 1 (notes placeholder
# "Espace réservé des notes 2"), reached via Slide.NotesPage.
$np = $s.NotesPage
$notesShape = $np.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

$firstPara  = "La première partie de ce projet a été réalisée par Sébastien Henneberger et Thibault Schowing (déjà notée)"
$secondPara = "La seconde partie, celle présentée dans ce Powerpoint, a été réalisée par Anastasia Zharkova et Thibault Schowing et consiste à sécuriser la partie 1."

# `n` starts a new paragraph; the trailing `n` keeps the final (empty)
# paragraph that was already present after the old TODO/CSRF note.
$tr.Text = $firstPara + "`n" + $secondPara + "`n"
